# chore: status log + auto-updated Current Price
$wb = $excel.ActiveWorkbook

# 1) Auto-updated Current Price on the Assumptions sheet (B5: 80.40000152587891 -> 79.15000152587891)
$wsAssump = $wb.Worksheets.Item("Assumptions")
$wsAssump.Range("B5").Value = 79.15000152587891

# 2) Append a new status snapshot row to the Status_Log sheet
$wsStatus = $wb.Worksheets.Item("Status_Log")

$wsStatus.Range("A4").Value = "2025-12-29T07:55:51Z"
$wsStatus.Range("B4").Value = 79.15000152587891
$wsStatus.Range("C4").Value = 15.31404656443955
$wsStatus.Range("D4").Value = 416.8457676605972
$wsStatus.Range("E4").Value = 0
$wsStatus.Range("F4").Value = 0
$wsStatus.Range("G4").Value = 0
$wsStatus.Range("H4").Value = 0
$wsStatus.Range("I4").Value = 0
$wsStatus.Range("J4").Value = "观察"
$wsStatus.Range("K4").Value = 0
